$d = $word.ActiveDocument

# 1) Fix the spelling of "Dillion" -> "Dillon" in the team lead name.
$d.Content.Find.Execute("Team Lead Dillion Johnson", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Team Lead Dillon Johnson", 2) | Out-Null

# 2) Remove the stray lastRenderedPageBreak on the "09/18/2018" run by
#    doing a self-replace on its text (forces the run to be rewritten
#    without the lastRenderedPageBreak child element).
$d.Content.Find.Execute("09/18/2018", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "09/18/2018", 2) | Out-Null

# 3) Append a second run of text to the "Slime interactions" bullet.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Slime interactions, plan to test each type of slime to see if it has the intended effects.") {
        $p.Range.InsertAfter(" How much slime can we fix of the screen at once before it crashes)")
        break
    }
}

# 4) Move the "_GoBack" bookmark from the very last paragraph of the
#    document up to the empty paragraph right after "QA Testing Lead
#    Carl Petersen" (this mirrors what Word does automatically when the
#    last edit position in the file changes).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "QA Testing Lead Carl Petersen") {
        $target = $paras.Item($i + 1)
        $d.Bookmarks.Add("_GoBack", $target.Range) | Out-Null
        break
    }
}
